$wb = $excel.ActiveWorkbook

# --- Metadata sheet updates ---
$meta = $wb.Worksheets.Item("Metadata")

# Version: 1.1.2 -> 1.1.3
$meta.Range("B3").Value = "1.1.3"

# Date: 2024-04-12T08:01:09+00:00 -> 2024-04-25T23:34:44+00:00
$meta.Range("B8").Value = "2024-04-25T23:34:44+00:00"

# Description: shorten to a single line
$meta.Range("B12").Value = "Logical Model for the HCERT"

# --- Elements sheet updates (HCert.5 / SMART Health Link row) ---
$elements = $wb.Worksheets.Item("Elements")

# Type(s) URL: ips-pilgrimage -> trust
$elements.Range("K6").Value = "http://smart.who.int/trust/StructureDefinition/SmartHealthLink`n"

# Definition: drop "(PROPOSED)"
$elements.Range("M6").Value = "SMART Health Link"
